# Atualizado por script em 20-11-2023 08:45
# Applies the betexplorer-style reshuffle of brazil serie-a 2023 rows and
# appends one newly scraped match (Fortaleza x Cruzeiro).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the match details (columns F:V) between a handful of row
#    pairs / cycles. The leading columns A:E (index, pais, torneio,
#    temporada, data_partida) stay put - only the home/away/odds/url
#    data moved around between rows.
# ---------------------------------------------------------------------

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("F$r1`:V$r1")
    $rng2 = $ws.Range("F$r2`:V$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Simple pairwise swaps
Swap-Rows 312 314
Swap-Rows 330 331
Swap-Rows 332 333

# 3-way rotation: 334 -> 335 -> 336 -> 334
$v334 = $ws.Range("F334:V334").Value2
$v335 = $ws.Range("F335:V335").Value2
$v336 = $ws.Range("F336:V336").Value2

$ws.Range("F335:V335").Value2 = $v334
$ws.Range("F336:V336").Value2 = $v335
$ws.Range("F334:V334").Value2 = $v336

# ---------------------------------------------------------------------
# 2) Append the new match row (337) at the bottom of the table.
# ---------------------------------------------------------------------

# Copy formatting (cell styles) from the last existing row so the new
# row matches the look of the rest of the table (bold index column,
# date-time formatted date column, ...).
$ws.Range("A336:V336").Copy() | Out-Null
$ws.Range("A337:V337").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A337").Value = 336
$ws.Range("B337").Value = "brazil"
$ws.Range("C337").Value = "serie-a"

$ws.Range("D337").Value = "'2023"
$ws.Range("D337").Style = "Normal"

$ws.Range("E337").Value = 45248.9375

$ws.Range("F337").Value = "Fortaleza"
$ws.Range("G337").Value = 0
$ws.Range("H337").Value = "Cruzeiro"
$ws.Range("I337").Value = 1

$ws.Range("J337").Value = 1.9
$ws.Range("K337").Value = "12/11/2023 22:42"
$ws.Range("L337").Value = 2.01
$ws.Range("M337").Value = "18/11/2023 22:21"

$ws.Range("N337").Value = 3.57
$ws.Range("O337").Value = "12/11/2023 22:42"
$ws.Range("P337").Value = 3.25
$ws.Range("Q337").Value = "18/11/2023 22:27"

$ws.Range("R337").Value = 4.28
$ws.Range("S337").Value = "12/11/2023 22:42"
$ws.Range("T337").Value = 4.5
$ws.Range("U337").Value = "18/11/2023 22:27"

$ws.Range("V337").Value = "https://www.betexplorer.com/football/brazil/serie-a/fortaleza-cruzeiro/riWa1h72/"

Write-Output "Row reshuffle + new row 337 applied"
